$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Spon2"
$ws.Range("C2").Value = "Itgb2"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.3355466666666667
$ws.Range("H2").Value = 1.00664
$ws.Range("I2").Value = 0.02388108083384315
$ws.Range("J2").Value = 0.02388108083384315
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 83.95844533333333
$ws.Range("N2").Value = 251.875336
$ws.Range("O2").Value = 0.9979754487867319
$ws.Range("P2").Value = 0.9979754487867319
$ws.Range("Q2").Value = 28.17197647011555
$ws.Range("R2").Value = 253.54778823104
$ws.Range("S2").Value = 0.02383273236266684
$ws.Range("T2").Value = 0.02383273236266684

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Spon2"
$ws.Range("C3").Value = "Itgb2"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.3355466666666667
$ws.Range("H3").Value = 1.00664
$ws.Range("I3").Value = 0.02388108083384315
$ws.Range("J3").Value = 0.02388108083384315
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.170323
$ws.Range("N3").Value = 0.510969
$ws.Range("O3").Value = 0.002024551213268089
$ws.Range("P3").Value = 0.00202455121326809
$ws.Range("Q3").Value = 0.05715131490666667
$ws.Range("R3").Value = 0.51436183416
$ws.Range("S3").Value = 0.00004834847117631046
$ws.Range("T3").Value = 0.00004834847117631048

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Spon2"
$ws.Range("C4").Value = "Itgb2"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 12.878362
$ws.Range("H4").Value = 38.635086
$ws.Range("I4").Value = 0.9165616424823987
$ws.Range("J4").Value = 0.9165616424823989
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 83.95844533333333
$ws.Range("N4").Value = 251.875336
$ws.Range("O4").Value = 0.9979754487867319
$ws.Range("P4").Value = 0.9979754487867319
$ws.Range("Q4").Value = 1081.247251959877
$ws.Range("R4").Value = 9731.225267638896
$ws.Range("S4").Value = 0.9147060164970759
$ws.Range("T4").Value = 0.9147060164970762

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Spon2"
$ws.Range("C5").Value = "Itgb2"
$ws.Range("D5").Value = "FAPs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 12.878362
$ws.Range("H5").Value = 38.635086
$ws.Range("I5").Value = 0.9165616424823987
$ws.Range("J5").Value = 0.9165616424823989
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 0.170323
$ws.Range("N5").Value = 0.510969
$ws.Range("O5").Value = 0.002024551213268089
$ws.Range("P5").Value = 0.00202455121326809
$ws.Range("Q5").Value = 2.193481250926
$ws.Range("R5").Value = 19.741331258334
$ws.Range("S5").Value = 0.001855625985322733
$ws.Range("T5").Value = 0.001855625985322734

# Row 6
$ws.Range("A6").Value = "sCs"
$ws.Range("B6").Value = "Spon2"
$ws.Range("C6").Value = "Itgb2"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.8368233333333333
$ws.Range("H6").Value = 2.51047
$ws.Range("I6").Value = 0.05955727668375805
$ws.Range("J6").Value = 0.05955727668375806
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 83.95844533333333
$ws.Range("N6").Value = 251.875336
$ws.Range("O6").Value = 0.9979754487867319
$ws.Range("P6").Value = 0.9979754487867319
$ws.Range("Q6").Value = 70.25838608532443
$ws.Range("R6").Value = 632.3254747679199
$ws.Range("S6").Value = 0.05943669992698901
$ws.Range("T6").Value = 0.05943669992698902

# Row 7
$ws.Range("A7").Value = "sCs"
$ws.Range("B7").Value = "Spon2"
$ws.Range("C7").Value = "Itgb2"
$ws.Range("D7").Value = "FAPs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.8368233333333333
$ws.Range("H7").Value = 2.51047
$ws.Range("I7").Value = 0.05955727668375805
$ws.Range("J7").Value = 0.05955727668375806
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.170323
$ws.Range("N7").Value = 0.510969
$ws.Range("O7").Value = 0.002024551213268089
$ws.Range("P7").Value = 0.00202455121326809
$ws.Range("Q7").Value = 0.1425302606033333
$ws.Range("R7").Value = 1.28277234543
$ws.Range("S7").Value = 0.0001205767567690456
$ws.Range("T7").Value = 0.0001205767567690457
